# Update "想去人数" (interest counts) in column F on both the "展览" and
# "全部类型" sheets, which mirror each other's data.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F3"  = 290
    "F4"  = 11076
    "F5"  = 10280
    "F12" = 24
    "F13" = 9598
    "F14" = 9
    "F17" = 37
    "F19" = 85
    "F20" = 390
    "F21" = 10858
    "F22" = 10782
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
